$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 6 (hunk 0)
$ws.Cells.Item(6, 8).Value = 2643.5833
$ws.Cells.Item(6, 10).Value = 2736.7827
$ws.Cells.Item(6, 12).Value = 8210.348100000001
$ws.Cells.Item(6, 14).Value = -8434.348100000001
# Row 9 (hunk 1)
$ws.Cells.Item(9, 8).Value = 96
$ws.Cells.Item(9, 9).Value = 96
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 96
$ws.Cells.Item(9, 12).Value = 0
$ws.Cells.Item(9, 13).Value = 73
$ws.Cells.Item(9, 14).Value = ""
# Row 12 (hunk 2)
$ws.Cells.Item(12, 8).Value = 1059.7
$ws.Cells.Item(12, 9).Value = 779.4
$ws.Cells.Item(12, 10).Value = 1340
$ws.Cells.Item(12, 11).Value = 779.4
$ws.Cells.Item(12, 12).Value = 1340
$ws.Cells.Item(12, 13).Value = -609.4
$ws.Cells.Item(12, 14).Value = -1680
# Row 21 (hunk 3)
$ws.Cells.Item(21, 8).Value = 43407.6
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 12).Value = 0
$ws.Cells.Item(21, 14).Value = ""
# Row 23 (hunk 4)
$ws.Cells.Item(23, 8).Value = 43407.6
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 12).Value = 0
$ws.Cells.Item(23, 14).Value = ""
# Row 38 (hunk 5)
$ws.Cells.Item(38, 8).Value = 1459.3334
$ws.Cells.Item(38, 9).Value = 1189
$ws.Cells.Item(38, 11).Value = 3567
$ws.Cells.Item(38, 13).Value = -3195
# Row 58 (hunk 6)
$ws.Cells.Item(58, 8).Value = 519.25
$ws.Cells.Item(58, 9).Value = 519.25
$ws.Cells.Item(58, 10).Value = 0
$ws.Cells.Item(58, 11).Value = 1557.75
$ws.Cells.Item(58, 12).Value = 0
$ws.Cells.Item(58, 13).Value = -1407.75
$ws.Cells.Item(58, 14).Value = ""
# Row 138 (hunk 7)
$ws.Cells.Item(138, 8).Value = 3315.9592
$ws.Cells.Item(138, 10).Value = 3253.1667
$ws.Cells.Item(138, 12).Value = 9759.500100000001
$ws.Cells.Item(138, 14).Value = -20039.5001
# Row 139 (hunk 8)
$ws.Cells.Item(139, 8).Value = 73737.375
$ws.Cells.Item(139, 10).Value = 73737.375
$ws.Cells.Item(139, 12).Value = 73737.375
$ws.Cells.Item(139, 14).Value = -84017.375
# Row 140 (hunk 9)
$ws.Cells.Item(140, 8).Value = 82399.39999999999
$ws.Cells.Item(140, 10).Value = 82399.39999999999
$ws.Cells.Item(140, 12).Value = 82399.39999999999
$ws.Cells.Item(140, 14).Value = -92759.39999999999
# Row 141 (hunk 10)
$ws.Cells.Item(141, 8).Value = 3502611.8
$ws.Cells.Item(141, 10).Value = 4999.5
$ws.Cells.Item(141, 12).Value = 14998.5
$ws.Cells.Item(141, 14).Value = -25358.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 4 (hunk 11)
$ws.Cells.Item(4, 8).Value = 367
$ws.Cells.Item(4, 9).Value = 299.5
$ws.Cells.Item(4, 11).Value = 299.5
$ws.Cells.Item(4, 13).Value = -183.5
# Row 32 (hunk 12)
$ws.Cells.Item(32, 8).Value = 15270.23
$ws.Cells.Item(32, 9).Value = 12418.4
$ws.Cells.Item(32, 10).Value = 24776.334
$ws.Cells.Item(32, 11).Value = 12418.4
$ws.Cells.Item(32, 12).Value = 24776.334
$ws.Cells.Item(32, 13).Value = -12131.4
$ws.Cells.Item(32, 14).Value = -25350.334
# Row 61 (hunk 13)
$ws.Cells.Item(61, 8).Value = 18658.438
$ws.Cells.Item(61, 9).Value = 29608.678
$ws.Cells.Item(61, 11).Value = 29608.678
$ws.Cells.Item(61, 13).Value = -29396.678
# Row 74 (hunk 14)
$ws.Cells.Item(74, 8).Value = 728.5
$ws.Cells.Item(74, 9).Value = 522.0909
$ws.Cells.Item(74, 11).Value = 522.0909
$ws.Cells.Item(74, 13).Value = 351.9091
# Row 77 (hunk 15)
$ws.Cells.Item(77, 8).Value = 728.5
$ws.Cells.Item(77, 9).Value = 522.0909
$ws.Cells.Item(77, 11).Value = 2610.4545
$ws.Cells.Item(77, 13).Value = 1757.5455
# Row 97 (hunk 16)
$ws.Cells.Item(97, 8).Value = 938.8889
$ws.Cells.Item(97, 9).Value = 868.75
$ws.Cells.Item(97, 11).Value = 868.75
$ws.Cells.Item(97, 13).Value = -372.75
# Row 110 (hunk 17)
$ws.Cells.Item(110, 8).Value = 1463.375
$ws.Cells.Item(110, 9).Value = 1249
$ws.Cells.Item(110, 11).Value = 1249
$ws.Cells.Item(110, 13).Value = 796
# Row 132 (hunk 18)
$ws.Cells.Item(132, 8).Value = 1938.5146
$ws.Cells.Item(132, 9).Value = 1575.225
$ws.Cells.Item(132, 11).Value = 4725.674999999999
$ws.Cells.Item(132, 13).Value = -2195.674999999999
# Row 136 (hunk 19)
$ws.Cells.Item(136, 8).Value = 18658.438
$ws.Cells.Item(136, 9).Value = 29608.678
$ws.Cells.Item(136, 11).Value = 88826.034
$ws.Cells.Item(136, 13).Value = -86276.034

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 99 (hunk 20)
$ws.Cells.Item(99, 8).Value = 891.6667
$ws.Cells.Item(99, 9).Value = 891.6667
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 891.6667
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).Value = 606.3333
$ws.Cells.Item(99, 14).Value = ""
# Row 141 (hunk 21)
$ws.Cells.Item(141, 8).Value = 65757.60000000001
$ws.Cells.Item(141, 10).Value = 67510.86
$ws.Cells.Item(141, 12).Value = 67510.86
$ws.Cells.Item(141, 14).Value = -77870.86

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31 (hunk 22)
$ws.Cells.Item(31, 8).Value = 2934.1304
$ws.Cells.Item(31, 9).Value = 1990.1538
$ws.Cells.Item(31, 10).Value = 4161.3
$ws.Cells.Item(31, 11).Value = 1990.1538
$ws.Cells.Item(31, 12).Value = 4161.3
$ws.Cells.Item(31, 13).Value = -1695.1538
$ws.Cells.Item(31, 14).Value = -4751.3
# Row 34 (hunk 23)
$ws.Cells.Item(34, 8).Value = 2934.1304
$ws.Cells.Item(34, 9).Value = 1990.1538
$ws.Cells.Item(34, 10).Value = 4161.3
$ws.Cells.Item(34, 11).Value = 1990.1538
$ws.Cells.Item(34, 12).Value = 4161.3
$ws.Cells.Item(34, 13).Value = -1788.1538
$ws.Cells.Item(34, 14).Value = -4565.3
# Row 58 (hunk 24)
$ws.Cells.Item(58, 8).Value = 854085.4399999999
$ws.Cells.Item(58, 9).Value = 2718729
$ws.Cells.Item(58, 11).Value = 2718729
$ws.Cells.Item(58, 13).Value = -2718526
# Row 136 (hunk 25)
$ws.Cells.Item(136, 8).Value = 854085.4399999999
$ws.Cells.Item(136, 9).Value = 2718729
$ws.Cells.Item(136, 11).Value = 8156187
$ws.Cells.Item(136, 13).Value = -8153637

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4 (hunk 26)
$ws.Cells.Item(4, 8).Value = 2716994.8
$ws.Cells.Item(4, 9).Value = 5250000
$ws.Cells.Item(4, 10).Value = 1450492
$ws.Cells.Item(4, 11).Value = 15750000
$ws.Cells.Item(4, 12).Value = 4351476
$ws.Cells.Item(4, 13).Value = -15749888
$ws.Cells.Item(4, 14).Value = -4351700
# Row 92 (hunk 27)
$ws.Cells.Item(92, 8).Value = 499.5
$ws.Cells.Item(92, 10).Value = 500
$ws.Cells.Item(92, 12).Value = 1500
$ws.Cells.Item(92, 14).Value = -3996
# Row 113 (hunk 28)
$ws.Cells.Item(113, 8).Value = 45194.92
$ws.Cells.Item(113, 10).Value = 1285.6111
$ws.Cells.Item(113, 12).Value = 3856.8333
$ws.Cells.Item(113, 14).Value = -8196.8333

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80 (hunk 29)
$ws.Cells.Item(80, 8).Value = 2930
$ws.Cells.Item(80, 9).Value = 2999.0908
$ws.Cells.Item(80, 11).Value = 2999.0908
$ws.Cells.Item(80, 13).Value = -2001.0908
# Row 83 (hunk 30)
$ws.Cells.Item(83, 8).Value = 2930
$ws.Cells.Item(83, 9).Value = 2999.0908
$ws.Cells.Item(83, 11).Value = 14995.454
$ws.Cells.Item(83, 13).Value = -10003.454
# Row 123 (hunk 31)
$ws.Cells.Item(123, 8).Value = 32472.4
$ws.Cells.Item(123, 10).Value = 32472.4
$ws.Cells.Item(123, 12).Value = 32472.4
$ws.Cells.Item(123, 14).Value = -37372.4
# Row 132 (hunk 32)
$ws.Cells.Item(132, 8).Value = 877916.9
$ws.Cells.Item(132, 9).Value = 990024.5
$ws.Cells.Item(132, 11).Value = 2970073.5
$ws.Cells.Item(132, 13).Value = -2967543.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81 (hunk 33)
$ws.Cells.Item(81, 8).Value = 699.6
$ws.Cells.Item(81, 9).Value = 699.6
$ws.Cells.Item(81, 11).Value = 1399.2
$ws.Cells.Item(81, 13).Value = -338.2
# Row 84 (hunk 34)
$ws.Cells.Item(84, 8).Value = 699.6
$ws.Cells.Item(84, 9).Value = 699.6
$ws.Cells.Item(84, 11).Value = 6996
$ws.Cells.Item(84, 13).Value = -1692

Write-Output "Edits applied"